# SeniorConnect_MasterLog.xlsx auto-update
# Appends the latest sensor/alert readings (2026-02-01, ~18:36-18:37) to the
# ALERTS, PIR, Humidity and Temperature logs.

$wb = $excel.ActiveWorkbook

function Set-LogRow {
    param(
        $ws,
        [int]$r,
        [string]$date,
        [string]$timestamp,
        [string]$hour,
        [string]$location,
        [string]$value,
        [string]$status
    )
    # Columns A-C (date/timestamp/hour) and any numeric-looking E values must
    # be forced to literal text -- otherwise Excel's COM layer auto-coerces
    # them into date/time/number types, same as typing them in the UI would
    # without a leading apostrophe.
    $ws.Cells.Item($r, 1).Value = "'" + $date
    $ws.Cells.Item($r, 2).Value = "'" + $timestamp
    $ws.Cells.Item($r, 3).Value = "'" + $hour
    $ws.Cells.Item($r, 4).Value = $location
    if ($value -match '^-?[0-9.,]+%?$') {
        $ws.Cells.Item($r, 5).Value = "'" + $value
    } else {
        $ws.Cells.Item($r, 5).Value = $value
    }
    $ws.Cells.Item($r, 6).Value = $status
}

# ---------------------------------------------------------------------------
# ALERTS sheet: two new alert rows (19-20)
# ---------------------------------------------------------------------------
$wsAlerts = $wb.Worksheets.Item("ALERTS")

Set-LogRow $wsAlerts 19 "2026-02-01" "18:36:34" "18:00" "Bathroom" "MODERATE" "MODERATE ALERT: Bathroom occupied, no motion > 40s."
Set-LogRow $wsAlerts 20 "2026-02-01" "18:36:54" "18:00" "Bathroom" "CRITICAL" "CRITICAL ALERT: Bathroom occupied, no motion > 60s."

# ---------------------------------------------------------------------------
# PIR sheet: thirteen new "No Motion" / "Inactive" rows (94-106)
# ---------------------------------------------------------------------------
$wsPir = $wb.Worksheets.Item("PIR")

$pirTimes = @(
    "18:36:23", "18:36:25", "18:36:28", "18:36:33", "18:36:38", "18:36:43",
    "18:36:48", "18:36:53", "18:36:58", "18:37:03", "18:37:08", "18:37:13",
    "18:37:18"
)
$r = 94
foreach ($t in $pirTimes) {
    Set-LogRow $wsPir $r "2026-02-01" $t "18:00" "Bathroom" "No Motion" "Inactive"
    $r++
}

# ---------------------------------------------------------------------------
# Humidity sheet: eleven new "Active" readings (149-159)
# ---------------------------------------------------------------------------
$wsHumidity = $wb.Worksheets.Item("Humidity")

$humidityRows = @(
    @("18:36:24", "80.8%"),
    @("18:36:29", "79.3%"),
    @("18:36:34", "79.9%"),
    @("18:36:39", "79.9%"),
    @("18:36:44", "79.6%"),
    @("18:36:50", "79.1%"),
    @("18:36:55", "79.0%"),
    @("18:37:00", "77.9%"),
    @("18:37:05", "78.9%"),
    @("18:37:10", "77.9%"),
    @("18:37:20", "78.0%")
)
$r = 149
foreach ($row in $humidityRows) {
    Set-LogRow $wsHumidity $r "2026-02-01" $row[0] "18:00" "Bathroom" $row[1] "Active"
    $r++
}

# ---------------------------------------------------------------------------
# Temperature sheet: eleven new "Active" readings (149-159)
# ---------------------------------------------------------------------------
$wsTemperature = $wb.Worksheets.Item("Temperature")

$temperatureRows = @(
    @("18:36:24", "30.5C"),
    @("18:36:30", "30.4C"),
    @("18:36:35", "30.4C"),
    @("18:36:40", "30.4C"),
    @("18:36:45", "30.4C"),
    @("18:36:50", "30.3C"),
    @("18:36:56", "30.3C"),
    @("18:37:00", "30.2C"),
    @("18:37:05", "30.3C"),
    @("18:37:10", "30.2C"),
    @("18:37:20", "30.1C")
)
$r = 149
foreach ($row in $temperatureRows) {
    Set-LogRow $wsTemperature $r "2026-02-01" $row[0] "18:00" "Bathroom" $row[1] "Active"
    $r++
}
